$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# Leading "'" forces text entry for numeric-looking price strings (e.g. "1.001")
# so they stay text cells like the rest of the sheet instead of becoming numbers.

$ws.Range("D2").Value = "23.219.34"
$ws.Range("E2").Value = "  +0.71%  "
$ws.Range("D3").Value = "1.603.28"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").Value = "'1.001"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "'304.88"
$ws.Range("E6").Value = "  +0.83%  "
$ws.Range("D7").Value = "'0.3761"
$ws.Range("E7").Value = "  -0.37%  "
$ws.Range("D8").Value = "'53.01"
$ws.Range("E8").Value = "  +4.19%  "
$ws.Range("D9").Value = "'0.3602"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").Value = "'1.258"
$ws.Range("E10").Value = "  +0.78%  "
$ws.Range("D11").Value = "'1.001"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").Value = "'0.08134"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").Value = "'22.81"
$ws.Range("E13").Value = "  +2.05%  "
$ws.Range("D14").Value = "'6.596"
$ws.Range("E14").Value = "  +0.46%  "
$ws.Range("D15").Value = "'7.340"
$ws.Range("E15").Value = "  -0.25%  "
$ws.Range("D16").Value = "'0.00001242"
$ws.Range("E16").Value = "  +0.03%  "
$ws.Range("D17").Value = "1.602.73"
$ws.Range("E17").Value = "  +0.29%  "
$ws.Range("D18").Value = "'93.95"
$ws.Range("E18").Value = "  +1.79%  "
$ws.Range("D19").Value = "'0.06917"
$ws.Range("E19").Value = "  +0.96%  "
$ws.Range("D20").Value = "'18.14"
$ws.Range("E20").Value = "  +0.01%  "
$ws.Range("D21").Value = "'6.522"
$ws.Range("E22").Value = "  +0.24%  "
$ws.Range("D23").Value = "'12.88"
$ws.Range("D24").Value = "23.219.19"
$ws.Range("E24").Value = "  +0.69%  "
$ws.Range("D25").Value = "'3.070"
$ws.Range("E25").Value = "  +9.71%  "
$ws.Range("D26").Value = "'2.417"
$ws.Range("E26").Value = "  +1.79%  "
$ws.Range("D27").Value = "'21.16"
$ws.Range("E27").Value = "  +0.33%  "
$ws.Range("D28").Value = "'150.46"
$ws.Range("E28").Value = "  +0.94%  "
$ws.Range("D29").Value = "'5.261"
$ws.Range("E29").Value = "  +0.11%  "
$ws.Range("D30").Value = "'134.94"
$ws.Range("E30").Value = "  -0.02%  "
$ws.Range("D31").Value = "'2.407"
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("D32").Value = "'6.728"
$ws.Range("E32").Value = "  -0.27%  "
$ws.Range("D33").Value = "1.781.41"
$ws.Range("E33").Value = "  +0.47%  "
$ws.Range("E34").Value = "  -1.18%  "
$ws.Range("D35").Value = "'0.02762"
$ws.Range("E35").Value = "  +1.93%  "
$ws.Range("D36").Value = "'0.07393"
$ws.Range("E36").Value = "  -1.72%  "
$ws.Range("D37").Value = "'10.27"
$ws.Range("E37").Value = "  +0.87%  "
$ws.Range("D38").Value = "'0.2511"
$ws.Range("E38").Value = "  -0.13%  "
$ws.Range("D39").Value = "'6.108"
$ws.Range("E39").Value = "  -1.48%  "
$ws.Range("D40").Value = "'0.08746"
$ws.Range("E40").Value = "  -0.80%  "
$ws.Range("D41").Value = "'1.401"
$ws.Range("E41").Value = "  +3.29%  "
$ws.Range("D42").Value = "'0.7087"
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("D43").Value = "'12.41"
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").Value = "'15.76"
$ws.Range("E44").Value = "  +4.15%  "
$ws.Range("D45").Value = "'0.6512"
$ws.Range("E45").Value = "  -0.94%  "
$ws.Range("E46").Value = "  +2.11%  "
$ws.Range("E47").Value = "  +0.19%  "
$ws.Range("D48").Value = "'133.86"
$ws.Range("E48").Value = "  +1.45%  "
$ws.Range("D49").Value = "'0.07958"
$ws.Range("E49").Value = "  +0.43%  "
$ws.Range("D50").Value = "'1.195"
$ws.Range("E50").Value = "  -1.80%  "
$ws.Range("D51").Value = "'1.187"
$ws.Range("E51").Value = "  -3.28%  "
